$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet-view selection (was B29 -> now A16)
$ws.Range("A16").Select()

# Row 10: GB Yellowtail Flounder - add Target TAC value
$ws.Range("D10").Value() = 15003

# Row 11: SNE/MA Yellowtail Flounder - add Target TAC value
$ws.Range("D11").Value() = 2092

# Row 13: Pollock - add Percent of TAC and Target TAC values
$ws.Range("B13").Value() = 8220
$ws.Range("D13").Value() = 54800

# Row 15: "Redfish - 500s" row becomes the "Redfish" row (text now matches
# the existing "Redfish" shared string used in column F, which also drops
# the now-unused "Redfish - 500s" shared string), with new data values
$ws.Range("A15").Value() = "Redfish"
$ws.Range("B15").Value() = 625
$ws.Range("D15").Value() = 62500

# Row 17: White Hake - add Target TAC value
$ws.Range("D17").Value() = 54914.28571428571

# Row 19: American Plaice - add Target TAC value, and reformat the cell
# (number format + font) to match the other populated Target TAC cells
$ws.Range("D19").Value() = 56709.090909090912
$ws.Range("D17").Copy()
$ws.Range("D19").PasteSpecial(-4122)

# Row 25: Ocean Pout - add Target TAC value
$ws.Range("D25").Value() = 56090.909090909088
